# Applies the "Add budget outputs with UD penalty" edit:
#  - Summary: sets Unmet Demand Penalty to 0.01 and refreshes the dependent
#    NPV / Unmet Demand / Household Surplus totals.
#  - Costs and Revenues: refreshes Total Revenues / Total Operation Variable
#    Costs / Total Profits for the affected years.
#  - DG Dispatch: the diesel generator now dispatches to cover unmet demand
#    (capped at its 400 kW capacity) for day-types 0.0-0.2 (rows 2-10).
#  - Unmet Demand: residual unmet demand after the new DG dispatch for the
#    same rows (2-10).
#  - Household Surplus: refreshed surplus values for years 0-2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 0.01
$wsSummary.Range("B6").Value = 192639.3288767941
$wsSummary.Range("B8").Value = 17785260.94054101
$wsSummary.Range("B10").Value = 5908589.326457601

# ---------------------------------------------------------------------
# Costs and Revenues
# ---------------------------------------------------------------------
$wsCosts = $wb.Worksheets.Item("Costs and Revenues")

# Row 2: Total Revenues
$wsCosts.Range("B2").Value = 450601.9862943711
$wsCosts.Range("C2").Value = 450601.9862943711
$wsCosts.Range("D2").Value = 450601.9862943711
$wsCosts.Range("F2").Value = 125752.60296485
$wsCosts.Range("G2").Value = 125752.60296485
$wsCosts.Range("H2").Value = 125752.60296485
$wsCosts.Range("I2").Value = 125752.60296485
$wsCosts.Range("J2").Value = 125752.60296485
$wsCosts.Range("K2").Value = 125752.60296485
$wsCosts.Range("L2").Value = 125752.60296485
$wsCosts.Range("M2").Value = 125752.60296485
$wsCosts.Range("N2").Value = 125752.60296485
$wsCosts.Range("O2").Value = 125752.60296485
$wsCosts.Range("P2").Value = 125752.60296485

# Row 4: Total Operation Variable Costs
$wsCosts.Range("B4").Value = 430898.5094426759
$wsCosts.Range("C4").Value = 430898.5094426759
$wsCosts.Range("D4").Value = 430898.5094426759

# Row 6: Total Profits
$wsCosts.Range("B6").Value = -13968.25892698585
$wsCosts.Range("C6").Value = -13968.25892698585
$wsCosts.Range("D6").Value = -13968.25892698585
$wsCosts.Range("E6").Value = -76736.57798215213
$wsCosts.Range("F6").Value = 56363.42201784793
$wsCosts.Range("G6").Value = 56363.42201784793
$wsCosts.Range("H6").Value = 56363.42201784793
$wsCosts.Range("I6").Value = 56363.42201784793
$wsCosts.Range("J6").Value = 56363.42201784793
$wsCosts.Range("K6").Value = 56363.42201784793
$wsCosts.Range("L6").Value = 56363.42201784793
$wsCosts.Range("M6").Value = 56363.42201784793
$wsCosts.Range("N6").Value = 56363.42201784793
$wsCosts.Range("O6").Value = 56363.42201784793
$wsCosts.Range("P6").Value = 56363.42201784793

# ---------------------------------------------------------------------
# DG Dispatch / Unmet Demand: rows 2-10 (day-types 0.0, 0.1, 0.2 repeated
# three times over), columns B..Y (hours 0..23).
# DG Dispatch now supplies min(oldUnmetDemand, 400); Unmet Demand keeps
# the residual max(oldUnmetDemand-400, 0).
# ---------------------------------------------------------------------
$wsDg = $wb.Worksheets.Item("DG Dispatch")
$wsUd = $wb.Worksheets.Item("Unmet Demand")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")

$dgPattern = @(
  @(382.7338416634806,365.2728917710076,354.683041620683,381.9303700722618,400,400,339.4748021157671,210.4758895704059,11.94928935461252,0,0,0,0,0,0,9.990699214544804,149.8691179411497,209.0200695862453,223.0958495641314,251.3456529078365,327.7522584701349,349.240968717413,369.731100678469,386.2379386560536),
  @(166.5331836498673,172.7084989883157,147.4450655646388,157.6450804554009,145.0692123933839,137.3435171632106,112.2354442364965,89.39663285141508,0.7465913262578567,0,0,0,0,0,0,0,100.1578341526431,171.6831711038378,200.1647286948216,225.9413820809748,232.8005871494253,251.6949831609196,205.7729852034775,205.6826957773044),
  @(179.8319801819373,167.2468210986278,148.6154730182124,146.4339626465692,145.4210480229312,167.9909793584588,162.2271725074396,155.4504749272583,93.35918011667277,22.26949182588285,0,0,0,0,2.721440735106512,86.16204325169439,177.2933913771695,224.0165980369723,227.9455894282815,286.3190293564909,252.137643323828,286.522998336591,225.7096553890372,218.5846533520948)
)

$udPattern = @(
  @(0,0,0,0,6.876045741711437,15.30273751513505,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
  @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
  @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
)

for ($blockStart = 2; $blockStart -le 8; $blockStart += 3) {
    for ($p = 0; $p -lt 3; $p++) {
        $r = $blockStart + $p
        for ($c = 0; $c -lt $cols.Length; $c++) {
            $colLetter = $cols[$c]
            $wsDg.Range("$colLetter$r").Value = $dgPattern[$p][$c]
            $wsUd.Range("$colLetter$r").Value = $udPattern[$p][$c]
        }
    }
}

# ---------------------------------------------------------------------
# Household Surplus: years 0-2
# ---------------------------------------------------------------------
$wsHh = $wb.Worksheets.Item("Household Surplus")
$wsHh.Range("B2").Value = 1069137.530885604
$wsHh.Range("B3").Value = 1069137.530885604
$wsHh.Range("B4").Value = 1069137.530885604

Write-Host "Edit applied."
